$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Solution Number"
$ws.Range("B1").Value = "Flags"
$ws.Range("C1").Value = "Mines"
$ws.Range("D1").Value = "Missing Mines"
$ws.Range("E1").Value = "Fitness"
$ws.Range("F1").Value = "# Flags in solution"
$ws.Range("G1").Value = "Total mines"
$ws.Range("H1").Value = "# Correct Flags"
$ws.Range("I1").Value = "% Correctly Identified"
$ws.Range("J1").Value = "Max Possible Fitness"

# ---------------------------------------------------------------------------
# Data rows (rows 2-6)
# Columns C, D, E, F, G, H hold text values (even the numeric-looking ones),
# so a leading apostrophe is used to force Excel to store them as text
# instead of auto-converting to numbers.
# ---------------------------------------------------------------------------

# Row 2 - Solution 1
$ws.Range("A2").Value = "Solution 1"
$ws.Range("B2").Value = "{(5, 5), (3, 4), (1, 5), (3, 1), (4, 6), (8, 6), (6, 3), (1, 3), (3, 5)}"
$ws.Range("C2").Value = "{(5, 5), (7, 1), (3, 4), (1, 5), (3, 1), (4, 6), (8, 6), (6, 3), (1, 3), (3, 5)}"
$ws.Range("D2").Value = "'(7, 1)"
$ws.Range("E2").Value = "'516"
$ws.Range("F2").Value = "'9"
$ws.Range("G2").Value = "'10"
$ws.Range("H2").Value = "'9"
$ws.Range("I2").Value = 90
$ws.Range("J2").Value = 571

# Row 3 - Solution 2
$ws.Range("A3").Value = "Solution 2"
$ws.Range("B3").Value = "{(5, 5), (7, 1), (1, 5), (3, 1), (4, 6), (8, 6), (6, 3), (1, 3), (3, 5)}"
$ws.Range("C3").Value = "{(5, 5), (7, 1), (3, 4), (1, 5), (3, 1), (4, 6), (8, 6), (6, 3), (1, 3), (3, 5)}"
$ws.Range("D3").Value = "'(3, 4)"
$ws.Range("E3").Value = "'516"
$ws.Range("F3").Value = "'9"
$ws.Range("G3").Value = "'10"
$ws.Range("H3").Value = "'9"
$ws.Range("I3").Value = 90
$ws.Range("J3").Value = 571

# Row 4 - Solution 3
$ws.Range("A4").Value = "Solution 3"
$ws.Range("B4").Value = "{(5, 5), (7, 1), (3, 4), (1, 5), (3, 1), (4, 6), (8, 6), (6, 3), (1, 3), (3, 5)}"
$ws.Range("C4").Value = "{(5, 5), (7, 1), (3, 4), (1, 5), (3, 1), (4, 6), (8, 6), (6, 3), (1, 3), (3, 5)}"
$ws.Range("D4").Value = "'"
$ws.Range("E4").Value = "'571"
$ws.Range("F4").Value = "'10"
$ws.Range("G4").Value = "'10"
$ws.Range("H4").Value = "'10"
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 571

# Row 5 - Solution 4
$ws.Range("A5").Value = "Solution 4"
$ws.Range("B5").Value = "{(5, 5), (7, 1), (3, 4), (1, 5), (3, 1), (4, 6), (6, 3), (1, 3), (3, 5)}"
$ws.Range("C5").Value = "{(5, 5), (7, 1), (3, 4), (1, 5), (3, 1), (4, 6), (8, 6), (6, 3), (1, 3), (3, 5)}"
$ws.Range("D5").Value = "'(8, 6)"
$ws.Range("E5").Value = "'516"
$ws.Range("F5").Value = "'9"
$ws.Range("G5").Value = "'10"
$ws.Range("H5").Value = "'9"
$ws.Range("I5").Value = 90
$ws.Range("J5").Value = 571

# Row 6 - Solution 5
$ws.Range("A6").Value = "Solution 5"
$ws.Range("B6").Value = "{(5, 5), (3, 4), (1, 5), (3, 1), (4, 6), (8, 6), (6, 3), (1, 3), (3, 5)}"
$ws.Range("C6").Value = "{(5, 5), (7, 1), (3, 4), (1, 5), (3, 1), (4, 6), (8, 6), (6, 3), (1, 3), (3, 5)}"
$ws.Range("D6").Value = "'(7, 1)"
$ws.Range("E6").Value = "'516"
$ws.Range("F6").Value = "'9"
$ws.Range("G6").Value = "'10"
$ws.Range("H6").Value = "'9"
$ws.Range("I6").Value = 90
$ws.Range("J6").Value = 571
